# Automatic update of files.
# - Bumps Taxonsorteringsordning (column B) from 79243 -> 79244 on all
#   "Garnlav" rows.
# - Re-syncs a handful of observation rows (12/13 and 15/16/17) whose
#   underlying records were renumbered/reordered by the source system,
#   carrying each record's full set of fields (Id, coordinates, accuracy,
#   start/end time, observer) along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column B: Taxonsorteringsordning 79243 -> 79244 --------------------
$bRows = @(2,3,4,5,6,7,9,12,13,15,16,17,19,20)
foreach ($r in $bRows) {
    $ws.Cells.Item($r, 2).Value = 79244
}

# ---- Rows 12 & 13: records swap places -----------------------------------
$ws.Range("A12").Value = 130979103
$ws.Range("Q12").Value = 570739
$ws.Range("R12").Value = 6736418
$ws.Range("S12").Value = 1
$ws.Range("Z12").Value = ""
$ws.Range("AB12").Value = ""
$ws.Range("AF12").Value = ""
$ws.Range("AW12").Value = "Erik Danielsson"
$ws.Range("AX12").Value = "Erik Danielsson"

$ws.Range("A13").Value = 130983072
$ws.Range("Q13").Value = 570809
$ws.Range("R13").Value = 6736404
$ws.Range("S13").Value = 10
$ws.Range("Z13").Value = "08:44"
$ws.Range("AB13").Value = "08:44"
$ws.Range("AF13").Value = ""
$ws.Range("AW13").Value = "Bo karlstens"
$ws.Range("AX13").Value = "Bo karlstens"

# ---- Rows 15, 16 & 17: records rotate (15<-17, 16<-15, 17<-16) ----------
$ws.Range("A15").Value = 130983071
$ws.Range("P15").Value = "Flytjärnsmyran, Dlr"
$ws.Range("Q15").Value = 570817
$ws.Range("R15").Value = 6736417
$ws.Range("Z15").Value = "08:53"
$ws.Range("AB15").Value = "08:53"
$ws.Range("AF15").Value = ""
$ws.Range("AW15").Value = "Bo karlstens"
$ws.Range("AX15").Value = "Bo karlstens"

$ws.Range("A16").Value = 130983619
$ws.Range("P16").Value = "Flytjärnsmyren, Dlr"
$ws.Range("Q16").Value = 570825
$ws.Range("R16").Value = 6736389
$ws.Range("Z16").Value = "08:54"
$ws.Range("AB16").Value = "08:54"
$ws.Range("AF16").Value = ""
$ws.Range("AW16").Value = "Göran Ehn"
$ws.Range("AX16").Value = "Göran Ehn"

$ws.Range("A17").Value = 130983074
$ws.Range("Q17").Value = 570764
$ws.Range("R17").Value = 6736425
$ws.Range("Z17").Value = "08:23"
$ws.Range("AB17").Value = "08:23"
$ws.Range("AW17").Value = "Bo karlstens"
$ws.Range("AX17").Value = "Bo karlstens"
